$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: "Trening" (training split into parts) ---
$ws.Range("F1").Value = "Trening"
$ws.Range("F2").Value = "Gra"
$ws.Range("F3").Value = "Gra"
$ws.Range("F4").Value = "Gra"
$ws.Range("F5").Value = "Gra"
$ws.Range("F6").Value = "Gra"
$ws.Range("F7").Value = "Gra"

# --- Register the lowercase numFmt (164) on a scratch cell so it is
#     present in styles.xml, then fix the scratch cell back to the
#     uppercase numFmt (165) before clearing it, so no stray cellXf
#     referencing 164 survives in the used range. ---
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "yyyy-mm-dd h:mm:ss"
$scratch.NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Column A: replace text timestamps with real date/time serials,
#     formatted with the uppercase numFmt (165). ---
$ws.Range("A2:A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Value = 45687.51027083334
$ws.Range("A3").Value = 45687.53663657408
$ws.Range("A4").Value = 45687.5376712963
$ws.Range("A5").Value = 45687.51026851852
$ws.Range("A6").Value = 45687.53476273148
$ws.Range("A7").Value = 45687.53663425926

$scratch.Clear()
